$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 38, which pushes existing rows 38-83 down to 39-84
$ws.Rows.Item(38).Insert()

# Populate the new row 38 with the new data record
$ws.Cells.Item(38, 1).Value = 5
$ws.Cells.Item(38, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(38, 3).Value = "Maule"
$ws.Cells.Item(38, 4).Value = 44546
$ws.Cells.Item(38, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38, 5).Value = 7
$ws.Cells.Item(38, 6).Value = 100112022
$ws.Cells.Item(38, 7).Value = "Arveja Verde"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 300
$ws.Cells.Item(38, 11).Value = 16000
$ws.Cells.Item(38, 12).Value = 16000
$ws.Cells.Item(38, 13).Value = 16000
$ws.Cells.Item(38, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Carahue"
$ws.Cells.Item(38, 16).Value = 640
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"
